$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value = 72
$ws.Range("R3").Value = 44
$ws.Range("Q10").Value = 72
$ws.Range("R10").Value = 40
$ws.Range("Q17").Value = 79
$ws.Range("R17").Value = 75
$ws.Range("Q23").Value = 6
$ws.Range("R23").Value = 3
$ws.Range("Q32").Value = 54
$ws.Range("R32").Value = 27
$ws.Range("Q40").Value = 18
$ws.Range("R40").Value = 15
$ws.Range("Q49").Value = 90
$ws.Range("R49").Value = 73
$ws.Range("Q58").Value = 57
$ws.Range("R58").Value = 57
$ws.Range("Q66").Value = 61
$ws.Range("R66").Value = 8
$ws.Range("Q74").Value = 86
$ws.Range("R74").Value = 30
$ws.Range("Q78").Value = 55
$ws.Range("R78").Value = 48
$ws.Range("Q89").Value = 100
$ws.Range("R89").Value = 97
$ws.Range("Q97").Value = 49
$ws.Range("R97").Value = 41
$ws.Range("Q106").Value = 96
$ws.Range("R106").Value = 7
$ws.Range("Q115").Value = 27
$ws.Range("R115").Value = 27
$ws.Range("Q124").Value = 76
$ws.Range("R124").Value = 14
$ws.Range("Q133").Value = 25
$ws.Range("R133").Value = 18
$ws.Range("Q142").Value = 65
$ws.Range("R142").Value = 60
